# No_conformidades.xlsx - "Generación de calidad de P1356"
#
# The non-conformity entry "Notificación linea base" / "Enviar notificación
# de creación de línea base" (Excel row 7) is resolved/removed from the
# tracker, so the whole sheet row is deleted. Excel shifts every following
# row up by one, renumbering the "ID" column (col A) and the sheet's
# dimension/validation ranges along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 7 ("Notificación linea base" non-conformity).
# This shifts rows 8:28 up to 7:27, carries row formatting/heights with
# them, renumbers the literal ID values in column A, shrinks the used
# range to A1:G27, and drops the two shared strings that were only
# referenced by that row.
$ws.Rows("7").Delete() | Out-Null

# Leave the selection where Excel lands after the delete.
$ws.Range("B5").Select() | Out-Null
